# Getting ready for PUGSUG demo.
# Update the title on slide 1 from "Git for VDW Work" to "Git for SAS Work".
#
# The original title is a single run: "Git for VDW Work" (16 chars).
#   "Git "     -> chars 1-4   (kept as-is)
#   "for VDW " -> chars 5-12  (retyped to "for SAS ")
#   "Work"     -> chars 13-16 (kept as-is)
#
# Re-typing just the middle chunk naturally leaves the untouched head/tail
# text in their own (split-off) runs while the retyped chunk becomes its
# own new run - matching how PowerPoint itself splits a run on a partial
# in-place edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$titleRange = $shape.TextFrame.TextRange

$middle = $titleRange.Characters(5, 8)
$middle.Text = "for SAS "
